$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price column cells so numeric-looking strings
# (e.g. "1.00", "685.20", multi-dot big numbers) stay exact text, matching
# the original inlineStr storage instead of being parsed into floats.
$priceCells = @("D2","D3","D5","D6","D7","D9","D11","D14","D15","D16","D17","D19","D20","D21","D23","D24","D25","D28","D31","D34","D35","D36","D37","D38","D39","D40","D41","D45","D46","D47","D48","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.276.53"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "3.681.71"
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "685.20"
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("D6").Value = "162.65"
$ws.Range("E6").Value = "  -5.63%  "
$ws.Range("D7").Value = "3.680.43"
$ws.Range("E7").Value = "  -3.72%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.499"
$ws.Range("D11").Value = "7.37"
$ws.Range("E11").Value = "  -3.71%  "
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("E13").Value = "  -5.80%  "
$ws.Range("D14").Value = "33.66"
$ws.Range("E14").Value = "  -6.69%  "
$ws.Range("D15").Value = "4.302.33"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "3.683.64"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").Value = "69.299.06"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "16.40"
$ws.Range("E19").Value = "  -6.43%  "
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  -7.70%  "
$ws.Range("D21").Value = "484.30"
$ws.Range("E22").Value = "  -7.16%  "
$ws.Range("D23").Value = "0.666"
$ws.Range("E23").Value = "  -8.18%  "
$ws.Range("D24").Value = "80.33"
$ws.Range("E24").Value = "  -5.07%  "
$ws.Range("D25").Value = "3.826.88"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("E26").Value = "  -9.98%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "11.43"
$ws.Range("E28").Value = "  -5.19%  "
$ws.Range("E29").Value = "  -9.01%  "
$ws.Range("E30").Value = "  -10.99%  "
$ws.Range("D31").Value = "2.72"
$ws.Range("E31").Value = "  -10.86%  "
$ws.Range("E32").Value = "  -7.82%  "
$ws.Range("E33").Value = "  -7.98%  "
$ws.Range("D34").Value = "27.17"
$ws.Range("E34").Value = "  -6.99%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.167"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").Value = "3.643.16"
$ws.Range("E37").Value = "  -3.94%  "
$ws.Range("D38").Value = "8.49"
$ws.Range("E38").Value = "  -7.70%  "
$ws.Range("D39").Value = "6.39"
$ws.Range("E39").Value = "  +6.73%  "
$ws.Range("D40").Value = "2.34"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "0.0933"
$ws.Range("E41").Value = "  -8.33%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -7.46%  "
$ws.Range("D45").Value = "163.30"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "48.39"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "2.83"
$ws.Range("E47").Value = "  -14.17%  "
$ws.Range("D48").Value = "29.76"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").Value = "0.000287"
$ws.Range("E49").Value = "  -8.59%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  -2.60%  "
